# Updates cryptos list values (Price / Volume(1h), plus a 3-row reorder
# for rows 42-44) to match the "Sun Sep 17 23:09:35 UTC 2023" GitHub Actions
# refresh of the crypto data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Text)
    $range = $ws.Range($Cell)
    # Some of the new values (e.g. "218.08", "2.32") parse as plain
    # numbers; force the cell to text first so Excel keeps them as
    # strings (matching every other cell in these columns), then drop
    # the number-format override again so no stray style sticks around.
    if ($Text -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $range.NumberFormat = "@"
        $range.Value = $Text
        $range.Style = "Normal"
    } else {
        $range.Value = $Text
    }
}

Set-TextValue 'D2' '26.704.28'
Set-TextValue 'E2' '  -0.31%  '

Set-TextValue 'D3' '1.633.07'
Set-TextValue 'E3' '  -1.00%  '

Set-TextValue 'E4' '  +0.03%  '

Set-TextValue 'D5' '218.08'
Set-TextValue 'E5' '  +0.56%  '

Set-TextValue 'E6' '  -1.61%  '

Set-TextValue 'E7' '  -0.01%  '

Set-TextValue 'E8' '  -1.40%  '

Set-TextValue 'E9' '  -1.18%  '

Set-TextValue 'D10' '18.94'
Set-TextValue 'E10' '  -1.73%  '

Set-TextValue 'D11' '0.0842'
Set-TextValue 'E11' '  -0.31%  '

Set-TextValue 'D12' '1.861.29'
Set-TextValue 'E12' '  -0.90%  '

Set-TextValue 'D13' '1.634.44'
Set-TextValue 'E13' '  -1.08%  '

Set-TextValue 'E14' '  -2.50%  '

Set-TextValue 'D15' '0.521'
Set-TextValue 'E15' '  -2.18%  '

Set-TextValue 'E16' '  -2.57%  '

Set-TextValue 'D17' '26.685.77'
Set-TextValue 'E17' '  -0.38%  '

Set-TextValue 'D18' '0.0₃0721'
Set-TextValue 'E18' '  -3.16%  '

Set-TextValue 'E19' '  +0.08%  '

Set-TextValue 'D20' '211.02'
Set-TextValue 'E20' '  -3.05%  '

Set-TextValue 'E21' '  -1.74%  '

Set-TextValue 'E22' '  -2.30%  '

Set-TextValue 'D23' '2.32'
Set-TextValue 'E23' '  -5.93%  '

Set-TextValue 'D24' '9.16'
Set-TextValue 'E24' '  -3.23%  '

Set-TextValue 'D25' '146.71'
Set-TextValue 'E25' '  +0.65%  '

Set-TextValue 'E26' '  +0.01%  '

Set-TextValue 'E27' '  -2.48%  '

Set-TextValue 'D28' '6.99'
Set-TextValue 'E28' '  -3.17%  '

Set-TextValue 'D29' '15.50'
Set-TextValue 'E29' '  -2.13%  '

Set-TextValue 'D30' '0.0500'
Set-TextValue 'E30' '  -4.19%  '

Set-TextValue 'E31' '  +0.83%  '

Set-TextValue 'E32' '  +0.20%  '

Set-TextValue 'E33' '  -2.91%  '

Set-TextValue 'D34' '1.259.55'
Set-TextValue 'E34' '  -1.43%  '

Set-TextValue 'D35' '1.52'
Set-TextValue 'E35' '  -2.35%  '

Set-TextValue 'E36' '  +0.26%  '

Set-TextValue 'E37' '  -3.48%  '

Set-TextValue 'E38' '  -3.52%  '

Set-TextValue 'E39' '  -0.06%  '

Set-TextValue 'E40' '  -4.06%  '

Set-TextValue 'D41' '0.798'
Set-TextValue 'E41' '  -2.67%  '

Set-TextValue 'B42' 'RocketPoolETH'
Set-TextValue 'C42' 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
Set-TextValue 'D42' '1.771.97'
Set-TextValue 'E42' '  -1.50%  '

Set-TextValue 'B43' 'MXToken'
Set-TextValue 'C43' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D43' '2.15'
Set-TextValue 'E43' '  -4.43%  '

Set-TextValue 'B44' 'FraxShare'
Set-TextValue 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue 'D44' '5.25'
Set-TextValue 'E44' '  -3.20%  '

Set-TextValue 'D45' '91.20'
Set-TextValue 'E45' '  -1.05%  '

Set-TextValue 'D46' '59.78'
Set-TextValue 'E46' '  +0.04%  '

Set-TextValue 'E47' '  -3.74%  '

Set-TextValue 'E48' '  -0.13%  '

Set-TextValue 'E49' '  -0.03%  '

Set-TextValue 'E50' '  -0.67%  '

Set-TextValue 'D51' '0.0954'
Set-TextValue 'E51' '  -2.79%  '
